$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.902.32"
$ws.Range("E2").Value = "  +0.97%  "

$ws.Range("D3").Value = "1.879.33"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("E4").Value = "  -0.53%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "324.94"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("E6").Value = "  -0.47%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4605"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3881"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.15%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07856"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.39%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.9859"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").Value = "1.889.34"
$ws.Range("E12").Value = "  -0.26%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "6.998"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -1.33%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.645"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.40%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.06966"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "87.99"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("E17").Value = "  -0.47%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000009986"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.70%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "16.99"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.50%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("D21").Value = "28.887.34"
$ws.Range("E21").Value = "  +0.89%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.233"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.13%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "10.96"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.85%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.086"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.41%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "156.25"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "19.34"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.36%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "6.033"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +2.40%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.927"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.90%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "117.48"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.09356"
$cell.NumberFormat = "General"
$cell.Style = "Normal"

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.9016"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -2.75%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.254"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.08%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.316"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.95%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.256"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.26%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.183"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +2.21%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.05738"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.02069"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("E38").Value = "  -0.54%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "7.641"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -4.59%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.5642"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.05%  "

$ws.Range("E41").Value = "  -1.94%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "9.656"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.51%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.272"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +4.39%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "11.82"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.5337"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.69%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.07046"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.47%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.841"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.25%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.531"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.81%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "112.50"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.060"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -5.21%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "70.61"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.61%  "
